$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("bass","drum","fx","guitar","little","piano","ride","string","synth")
$cols = @("I","J","K","L","M","N","O","P","Q")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "1").Value = $headers[$i]
}

$data = @(
    @(1,0,0,0,0,1,1,0,1),
    @(1,0,0,0,1,1,1,1,0),
    @(1,1,0,1,0,1,1,0,0),
    @(1,0,0,0,0,1,1,1,0),
    @(1,1,0,0,0,1,1,0,1),
    @(1,1,0,0,0,1,1,0,0),
    @(0,0,1,0,0,0,0,0,0),
    @(1,0,0,0,0,1,1,0,0)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $r + 2
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $ws.Range($cols[$c] + $row).Value = $data[$r][$c]
    }
}

$ws.Range("M1:N1").Font.Bold = $false

# Update view state (scroll position / active selection) to match the saved workbook.
$win = $ws.Application.ActiveWindow
$win.ScrollColumn = 6
$win.ScrollRow = 1
$null = $ws.Range("M7").Select()
